$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh drops the oldest historical row (previous row 62) and shifts
# all subsequent rows up by one date-slot, inserting a brand-new latest reading in
# row 63 and pushing a duplicate of the old last row down into a new row 95.
$data = @{
    63 = @{ D=44596; J=200; K=10000; L=11000; M=10500; N='$/caja 50 unidades'; O='Región de O''Higgins'; P=210; Q=50 }
    64 = @{ D=44238; J=100; K=10000; L=11000; M=10500; N='$/caja 60 unidades'; O='Región de O''Higgins'; P=175; Q=60 }
    65 = @{ D=44166; J=200; K=6000; L=7000; M=6500; N='$/caja 50 unidades'; O='Región de O''Higgins'; P=130; Q=50 }
    66 = @{ D=44223; J=100; K=9000; L=10000; M=9500; N='$/caja 60 unidades'; O='Región de O''Higgins'; P=158; Q=60 }
    67 = @{ D=44566; J=150; K=5500; L=6500; M=5967; N='$/caja 50 unidades'; O='Región Metropolitana'; P=119; Q=50 }
    68 = @{ D=44237; J=200; K=8000; L=9000; M=8500; N='$/caja 60 unidades'; O='Región de O''Higgins'; P=142; Q=60 }
    69 = @{ D=44195; J=200; K=10000; L=11000; M=10500; N='$/caja 60 unidades'; O='Región de O''Higgins'; P=175; Q=60 }
    70 = @{ D=44281; J=100; K=9000; L=10000; M=9500; N='$/caja 60 unidades'; O='Región de O''Higgins'; P=158; Q=60 }
    71 = @{ D=44420; J=100; K=9000; L=10000; M=9500; N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=190; Q=50 }
    72 = @{ D=44162; J=200; K=6000; L=6500; M=6250; N='$/caja 60 unidades'; O='Región de O''Higgins'; P=104; Q=60 }
    73 = @{ D=44343; J=100; K=9000; L=10000; M=9500; N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=190; Q=50 }
    74 = @{ D=44516; J=350; K=4500; L=5000; M=4786; N='$/caja 60 unidades'; O='Región de Arica y Parinacota'; P=80; Q=60 }
    75 = @{ D=44336; J=100; K=10000; L=12000; M=11000; N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=220; Q=50 }
    76 = @{ D=44565; J=100; K=7000; L=8000; M=7500; N='$/caja 50 unidades'; O='Región de O''Higgins'; P=150; Q=50 }
    77 = @{ D=44334; J=100; K=11000; L=12000; M=11500; N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=230; Q=50 }
    78 = @{ D=44357; J=100; K=8000; L=9000; M=8500; N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=170; Q=50 }
    79 = @{ D=44455; J=100; K=16000; L=17000; M=16500; N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=330; Q=50 }
    80 = @{ D=44537; J=170; K=5500; L=6000; M=5765; N='$/caja 60 unidades'; O='Región Metropolitana'; P=96; Q=60 }
    81 = @{ D=44397; J=100; K=8000; L=9000; M=8500; N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=170; Q=50 }
    82 = @{ D=44329; J=900; K=350; L=12000; M=1633; N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=33; Q=50 }
    83 = @{ D=44208; J=100; K=11000; L=12000; M=11500; N='$/caja 60 unidades'; O='Región de O''Higgins'; P=192; Q=60 }
    84 = @{ D=44355; J=100; K=9000; L=10000; M=9500; N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=190; Q=50 }
    85 = @{ D=44530; J=350; K=5000; L=5500; M=5214; N='$/caja 60 unidades'; O='Región del Maule'; P=87; Q=60 }
    86 = @{ D=44489; J=100; K=8000; L=9000; M=8500; N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=170; Q=50 }
    87 = @{ D=44572; J=260; K=14000; L=15000; M=14538; N='$/caja 60 unidades'; O='Región de O''Higgins'; P=242; Q=60 }
    88 = @{ D=44370; J=100; K=10000; L=11000; M=10500; N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=210; Q=50 }
    89 = @{ D=44385; J=100; K=9000; L=10000; M=9500; N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=190; Q=50 }
    90 = @{ D=44272; J=100; K=9000; L=10000; M=9500; N='$/caja 60 unidades'; O='Región de O''Higgins'; P=158; Q=60 }
    91 = @{ D=44299; J=100; K=7000; L=8000; M=7500; N='$/caja 50 unidades'; O='Región Metropolitana'; P=150; Q=50 }
    92 = @{ D=44258; J=200; K=10000; L=11000; M=10500; N='$/caja 60 unidades'; O='Región de O''Higgins'; P=175; Q=60 }
    93 = @{ D=44390; J=100; K=9000; L=10000; M=9500; N='$/caja 50 unidades'; O='Región de Arica y Parinacota'; P=190; Q=50 }
    94 = @{ D=44285; J=100; K=9000; L=10000; M=9500; N='$/caja 60 unidades'; O='Región de O''Higgins'; P=158; Q=60 }
    95 = @{ D=44498; J=350; K=10000; L=11000; M=10571; N='$/caja 60 unidades'; O='Región de Arica y Parinacota'; P=176; Q=60 }
}

# Row 95 is new: every data row shares the same market/category metadata, so
# seed those constant columns directly (mirrors row 94's A,B,C,E,F,G,H,I,R).
$ws.Cells.Item(95, 1).Value = 11
$ws.Cells.Item(95, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(95, 3).Value = "Bíobío"
$ws.Cells.Item(95, 4).NumberFormat = $ws.Cells.Item(94, 4).NumberFormat
$ws.Cells.Item(95, 5).Value = 8
$ws.Cells.Item(95, 6).Value = 100112032
$ws.Cells.Item(95, 7).Value = "Zapallo italiano"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 18).Value = "Hortaliza"

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 10).Value = $vals.J
    $ws.Cells.Item($row, 11).Value = $vals.K
    $ws.Cells.Item($row, 12).Value = $vals.L
    $ws.Cells.Item($row, 13).Value = $vals.M
    $ws.Cells.Item($row, 14).Value = $vals.N
    $ws.Cells.Item($row, 15).Value = $vals.O
    $ws.Cells.Item($row, 16).Value = $vals.P
    $ws.Cells.Item($row, 17).Value = $vals.Q
}
